$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.25000000000051
$ws.Range("H2").Value = 0.000001895722738010086
$ws.Range("I2").Value = 0.000001895722738010086
$ws.Range("L2").Value = 59.66221048547899
$ws.Range("M2").Value = "[32.57714135249849, 86.74727961845949]"
$ws.Range("N2").Value = 0.00005842482304929675
$ws.Range("O2").Value = 0.00005842482304929675
$ws.Range("P2").Value = 1.075500187715424
$ws.Range("Q2").Value = "[0.6100790538502707, 1.5409213215805782]"
$ws.Range("R2").Value = 0.00002881253693542973
$ws.Range("S2").Value = 0.00002881253693542973
$ws.Range("T2").Value = 61.43154385530821
$ws.Range("U2").Value = "[47.17948452922998, 75.68360318138645]"
$ws.Range("V2").Value = 0.00000000003588795927100819
$ws.Range("W2").Value = 0.00000000003588795927100819
$ws.Range("X2").Value = 20.92792792792835
$ws.Range("Y2").Value = 19.05755755755795
$ws.Range("Z2").Value = 22.79829829829876
$ws.Range("F3").Value = 25.25000000000051
$ws.Range("H3").Value = 0.00000269810728081854
$ws.Range("I3").Value = 0.00000269810728081854
$ws.Range("L3").Value = 56.87253875958717
$ws.Range("M3").Value = "[30.32285830655057, 83.42221921262376]"
$ws.Range("N3").Value = 0.00008648873359717868
$ws.Range("O3").Value = 0.00008648873359717868
$ws.Range("P3").Value = 1.012605439895809
$ws.Range("Q3").Value = "[0.5471843060306538, 1.4780265737609648]"
$ws.Range("R3").Value = 0.00006964518546248577
$ws.Range("S3").Value = 0.00006964518546248577
$ws.Range("T3").Value = 61.89071627247141
$ws.Range("U3").Value = "[48.08469576113013, 75.6967367838127]"
$ws.Range("V3").Value = 0.00000000001153210860138643
$ws.Range("W3").Value = 0.00000000001153210860138643
$ws.Range("X3").Value = 21.18068068068111
$ws.Range("Y3").Value = 19.31031031031069
$ws.Range("Z3").Value = 23.05105105105152
$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 25.25000000000051
$ws.Range("H4").Value = 0.0003074600776005187
$ws.Range("I4").Value = 0.0003074600776005187
$ws.Range("L4").Value = 52.34620857141648
$ws.Range("M4").Value = "[19.69824734204157, 84.99416980079138]"
$ws.Range("N4").Value = 0.00231974209031649
$ws.Range("O4").Value = 0.00231974209031649
$ws.Range("P4").Value = 0.9119738433844242
$ws.Range("Q4").Value = "[0.3333421634439615, 1.490605523324887]"
$ws.Range("R4").Value = 0.002708331939585085
$ws.Range("S4").Value = 0.002708331939585085
$ws.Range("T4").Value = 62.37880955676815
$ws.Range("U4").Value = "[45.544985768181846, 79.21263334535446]"
$ws.Range("V4").Value = 0.000000002104898033650215
$ws.Range("W4").Value = 0.000000002104898033650215
$ws.Range("X4").Value = 21.58508508508552
$ws.Range("Y4").Value = 19.25975975976015
$ws.Range("Z4").Value = 23.91041041041089
$ws.Range("B5").Value = 0
$ws.Range("F5").Value = 25.25000000000051
$ws.Range("H5").Value = 0.008929503845696574
$ws.Range("I5").Value = 0.008929503845696574
$ws.Range("L5").Value = 35.94764490943226
$ws.Range("M5").Value = "[6.9837368884765, 64.91155293038801]"
$ws.Range("N5").Value = 0.01614227333938789
$ws.Range("O5").Value = 0.01614227333938789
$ws.Range("P5").Value = 0.5094474573388847
$ws.Range("Q5").Value = "[-0.3522105877898465, 1.371105502467616]"
$ws.Range("R5").Value = 0.2399678817421349
$ws.Range("S5").Value = 0.2399678817421349
$ws.Range("T5").Value = 66.7655245223981
$ws.Range("U5").Value = "[51.328027251878986, 82.20302179291721]"
$ws.Range("V5").Value = 0.00000000003259903458285862
$ws.Range("W5").Value = 0.00000000003259903458285862
$ws.Range("X5").Value = 23.20270270270317
$ws.Range("Y5").Value = 19.73998998999039
$ws.Range("Z5").Value = 26.66541541541595
$ws.Range("F6").Value = 25.25000000000051
$ws.Range("H6").Value = 0.004632973754676084
$ws.Range("I6").Value = 0.004632973754676084
$ws.Range("L6").Value = 46.83658882795964
$ws.Range("M6").Value = "[13.695006034393458, 79.97817162152582]"
$ws.Range("N6").Value = 0.006637254244109103
$ws.Range("O6").Value = 0.006637254244109103
$ws.Range("P6").Value = 0.2830263651882694
$ws.Range("Q6").Value = "[-0.5597632555945777, 1.1258159859711165]"
$ws.Range("R6").Value = 0.5022625480788236
$ws.Range("S6").Value = 0.5022625480788236
$ws.Range("T6").Value = 67.02943057522316
$ws.Range("U6").Value = "[48.399305958567126, 85.65955519187919]"
$ws.Range("V6").Value = 0.000000004396834230036006
$ws.Range("W6").Value = 0.000000004396834230036006
$ws.Range("X6").Value = 24.1126126126131
$ws.Range("Y6").Value = 20.72572572572614
$ws.Range("Z6").Value = 27.49949949950005
$ws.Range("F7").Value = 25.25000000000051
$ws.Range("H7").Value = 0.001239142248212755
$ws.Range("I7").Value = 0.001239142248212755
$ws.Range("L7").Value = 43.88045872423756
$ws.Range("M7").Value = "[13.735011009897008, 74.02590643857812]"
$ws.Range("N7").Value = 0.005281799303205048
$ws.Range("O7").Value = 0.005281799303205048
$ws.Range("P7").Value = -0.6037895790683088
$ws.Range("Q7").Value = "[-1.2830528555201557, 0.07547369738353815]"
$ws.Range("R7").Value = 0.08013264357139427
$ws.Range("S7").Value = 0.08013264357139427
$ws.Range("T7").Value = 64.65849197335655
$ws.Range("U7").Value = "[48.80140383023466, 80.51558011647843]"
$ws.Range("V7").Value = 0.0000000001693669648972218
$ws.Range("W7").Value = 0.0000000001693669648972218
$ws.Range("X7").Value = 2.426426426426477
$ws.Range("Y7").Value = -0.3033033033033066
$ws.Range("Z7").Value = 5.156156156156261
$ws.Range("F8").Value = 25.45000000000054
$ws.Range("H8").Value = 0.00004352541515761921
$ws.Range("I8").Value = 0.00004352541515761921
$ws.Range("L8").Value = 44.59758003159163
$ws.Range("M8").Value = "[19.89720004092658, 69.29796002225669]"
$ws.Range("N8").Value = 0.0007078111493774575
$ws.Range("O8").Value = 0.0007078111493774575
$ws.Range("P8").Value = -0.8931054190385392
$ws.Range("Q8").Value = "[-1.434000250287232, -0.3522105877898465]"
$ws.Range("R8").Value = 0.001762562262686895
$ws.Range("S8").Value = 0.001762562262686895
$ws.Range("T8").Value = 61.54150630277946
$ws.Range("U8").Value = "[48.64113037779408, 74.44188222776484]"
$ws.Range("V8").Value = 0.000000000001790789738720377
$ws.Range("W8").Value = 0.000000000001790789738720377
$ws.Range("X8").Value = 3.617517517517594
$ws.Range("Y8").Value = 1.426626626626659
$ws.Range("Z8").Value = 5.80840840840853
$ws.Range("F9").Value = 25.45000000000054
$ws.Range("H9").Value = 0.001378489277603756
$ws.Range("I9").Value = 0.001378489277603756
$ws.Range("L9").Value = 51.50293556281895
$ws.Range("M9").Value = "[20.696475227109318, 82.30939589852858]"
$ws.Range("N9").Value = 0.001563416426712649
$ws.Range("O9").Value = 0.001563416426712649
$ws.Range("P9").Value = -1.559789745926464
$ws.Range("Q9").Value = "[-2.364842518017542, -0.754736973835386]"
$ws.Range("R9").Value = 0.0003156624976066169
$ws.Range("S9").Value = 0.0003156624976066169
$ws.Range("T9").Value = 68.51033489453204
$ws.Range("U9").Value = "[49.210074899590246, 87.81059488947383]"
$ws.Range("V9").Value = 0.00000000611942341244287
$ws.Range("W9").Value = 0.00000000611942341244287
$ws.Range("X9").Value = 6.317917917918052
$ws.Range("Y9").Value = 3.057057057057122
$ws.Range("Z9").Value = 9.578778778778982
$ws.Range("F10").Value = 25.45000000000054
$ws.Range("H10").Value = 0.0007823125990416013
$ws.Range("I10").Value = 0.0007823125990416013
$ws.Range("L10").Value = 46.37107065694269
$ws.Range("M10").Value = "[20.072205468173067, 72.66993584571232]"
$ws.Range("N10").Value = 0.0009122142330559857
$ws.Range("O10").Value = 0.0009122142330559857
$ws.Range("P10").Value = -1.610105544182156
$ws.Range("Q10").Value = "[-2.364842518017542, -0.8553685703467702]"
$ws.Range("R10").Value = 0.00009151743758617492
$ws.Range("S10").Value = 0.00009151743758617492
$ws.Range("T10").Value = 74.58533945850336
$ws.Range("U10").Value = "[58.052145777409, 91.1185331395977]"
$ws.Range("V10").Value = 0.000000000009580336524095401
$ws.Range("W10").Value = 0.000000000009580336524095401
$ws.Range("X10").Value = 6.52172172172186
$ws.Range("Y10").Value = 3.464664664664737
$ws.Range("Z10").Value = 9.578778778778982
$ws.Range("F11").Value = 25.45000000000054
$ws.Range("H11").Value = 0.02674208059126804
$ws.Range("I11").Value = 0.02674208059126804
$ws.Range("L11").Value = 31.55072471419547
$ws.Range("M11").Value = "[3.1664165457976736, 59.935032882593276]"
$ws.Range("N11").Value = 0.03016030674582804
$ws.Range("O11").Value = 0.03016030674582804
$ws.Range("P11").Value = -1.660421342437848
$ws.Range("Q11").Value = "[-2.880579450138388, -0.44026323473730855]"
$ws.Range("R11").Value = 0.008758197640762688
$ws.Range("S11").Value = 0.008758197640762688
$ws.Range("T11").Value = 49.73028216478083
$ws.Range("U11").Value = "[33.30038823407597, 66.1601760954857]"
$ws.Range("V11").Value = 0.0000002246131707295262
$ws.Range("W11").Value = 0.0000002246131707295262
$ws.Range("X11").Value = 6.725525525525669
$ws.Range("Y11").Value = 1.783283283283324
$ws.Range("Z11").Value = 11.66776776776801
$ws.Range("F12").Value = 25.45000000000054
$ws.Range("H12").Value = 0.00007536757438875163
$ws.Range("I12").Value = 0.00007536757438875163
$ws.Range("L12").Value = 60.32980989851517
$ws.Range("M12").Value = "[29.1857221712758, 91.47389762575453]"
$ws.Range("N12").Value = 0.000316409492217673
$ws.Range("O12").Value = 0.000316409492217673
$ws.Range("P12").Value = -2.868000500574466
$ws.Range("Q12").Value = "[-3.3963163822592364, -2.339684618889695]"
$ws.Range("R12").Value = 0.00000000000002953193245502916
$ws.Range("S12").Value = 0.00000000000002953193245502916
$ws.Range("T12").Value = 69.72780291456066
$ws.Range("U12").Value = "[52.58167994346074, 86.87392588566057]"
$ws.Range("V12").Value = 0.0000000001822202388979122
$ws.Range("W12").Value = 0.0000000001822202388979122
$ws.Range("X12").Value = 11.61681681681706
$ws.Range("Y12").Value = 9.476876876877075
$ws.Range("Z12").Value = 13.75675675675705
$ws.Range("F13").Value = 25.45000000000054
$ws.Range("H13").Value = 0.002205591893964187
$ws.Range("I13").Value = 0.002205591893964187
$ws.Range("L13").Value = 49.52793517972903
$ws.Range("M13").Value = "[16.097135612962433, 82.95873474649564]"
$ws.Range("N13").Value = 0.004586319979941189
$ws.Range("O13").Value = 0.004586319979941189
$ws.Range("P13").Value = -2.893158399702312
$ws.Range("Q13").Value = "[-3.673053272665544, -2.11326352673908]"
$ws.Range("R13").Value = 0.000000002046593783333606
$ws.Range("S13").Value = 0.000000002046593783333606
$ws.Range("T13").Value = 67.66794891640424
$ws.Range("U13").Value = "[49.43605906311274, 85.89983876969575]"
$ws.Range("V13").Value = 0.000000002021063538748535
$ws.Range("W13").Value = 0.000000002021063538748535
$ws.Range("X13").Value = 11.71871871871897
$ws.Range("Y13").Value = 8.55975975975994
$ws.Range("Z13").Value = 14.87767767767799
$ws.Range("F14").Value = 25.45000000000054
$ws.Range("H14").Value = 0.0000171468098824068
$ws.Range("I14").Value = 0.0000171468098824068
$ws.Range("L14").Value = 66.57904783497762
$ws.Range("M14").Value = "[36.89182361599731, 96.26627205395792]"
$ws.Range("N14").Value = 0.00004504899726942568
$ws.Range("O14").Value = 0.00004504899726942568
$ws.Range("P14").Value = -2.993789996213696
$ws.Range("Q14").Value = "[-3.522105877898466, -2.4654741145289267]"
$ws.Range("R14").Value = 0.000000000000007105427357601002
$ws.Range("S14").Value = 0.000000000000007105427357601002
$ws.Range("T14").Value = 73.07673956457812
$ws.Range("U14").Value = "[55.78376776976356, 90.36971135939268]"
$ws.Range("V14").Value = 0.00000000006289746501408899
$ws.Range("W14").Value = 0.00000000006289746501408899
$ws.Range("X14").Value = 12.12632632632658
$ws.Range("Y14").Value = 9.986386386386599
$ws.Range("Z14").Value = 14.26626626626657
